# Artfynd export refresh.
#
# Column B ("Taxonsorteringsordning") is bumped by 1 for every record that
# currently carries the "Garnlav"/Alectoria sarmentosa sort key (79243),
# and a handful of rows swap places because the underlying records they
# describe were re-ordered upstream. Apply both effects cell-by-cell so the
# rest of each row's data (coordinates, comments, flags, ...) travels with
# its record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that only get the plain Garnlav sort-order bump (79243 -> 79244) ---
foreach ($r in 2,5,7,11,16,17,21,22,23,24,25,26,27) {
    $ws.Range("B$r").Value = 79244
}

# --- Rows 3 <-> 4 swap places (Garnlav record <-> Tallticka record) ---
$ws.Range("A3").Value = 131041641
$ws.Range("B3").Value = 79244
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("P3").Value = "Tandbergsvasseln, Dlr"
$ws.Range("Q3").Value = 479078
$ws.Range("R3").Value = 6791615
$ws.Range("S3").Value = 50
$ws.Range("AC3").Value = "Rikligt i en radie av ca 50 meter, synfältet"

$ws.Range("A4").Value = 131039759
$ws.Range("B4").Value = 91830
$ws.Range("E4").Value = 5442
$ws.Range("F4").Value = "Tallticka"
$ws.Range("G4").Value = "Porodaedalea pini"
$ws.Range("H4").Value = "(Brot.) Murrill"
$ws.Range("P4").Value = "Gotvad, Dlr"
$ws.Range("Q4").Value = 479059
$ws.Range("R4").Value = 6792254
$ws.Range("S4").Value = 10
$ws.Range("AC4").ClearContents()

# --- Rows 9 <-> 10 swap places (Garnlav record <-> Mörk kolflarnlav record) ---
$ws.Range("A9").Value = 131039523
$ws.Range("B9").Value = 79244
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 479079
$ws.Range("R9").Value = 6792517

$ws.Range("A10").Value = 131040374
$ws.Range("B10").Value = 79002
$ws.Range("E10").Value = 228912
$ws.Range("F10").Value = "Mörk kolflarnlav"
$ws.Range("G10").Value = "Carbonicola myrmecina"
$ws.Range("H10").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q10").Value = 479088
$ws.Range("R10").Value = 6792211

# --- Rows 12 <-> 13 swap places (Garnlav record <-> Blanksvart spiklav record) ---
$ws.Range("A12").Value = 131040483
$ws.Range("B12").Value = 78647
$ws.Range("E12").Value = 6437
$ws.Range("F12").Value = "Blanksvart spiklav"
$ws.Range("G12").Value = "Calicium denigratum"
$ws.Range("H12").Value = "(Vain.) Tibell"
$ws.Range("Q12").Value = 479088
$ws.Range("R12").Value = 6792211

$ws.Range("A13").Value = 131039119
$ws.Range("B13").Value = 79244
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("Q13").Value = 479105
$ws.Range("R13").Value = 6792638

# --- Rows 18 -> 19 -> 20 -> 18 rotate (two Tretåig hackspett records + one Garnlav record) ---
$ws.Range("A18").Value = 131039579
$ws.Range("B18").Value = 79244
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("M18").ClearContents()
$ws.Range("P18").Value = "Gotvad, Dlr"
$ws.Range("Q18").Value = 479079
$ws.Range("R18").Value = 6792475
$ws.Range("S18").Value = 50
$ws.Range("AC18").Value = "Rikligt till måttligt i en radie av ca 50 meter, synfältet"

$ws.Range("A19").Value = 131041965
$ws.Range("M19").Value = "färska spår"
$ws.Range("P19").Value = "Tandbergsvasseln, Dlr"
$ws.Range("Q19").Value = 479096
$ws.Range("R19").Value = 6792085
$ws.Range("AE19").Value = $false

$ws.Range("A20").Value = 131039828
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("M20").Value = "bobygge"
$ws.Range("Q20").Value = 479059
$ws.Range("R20").Value = 6792254
$ws.Range("S20").Value = 10
$ws.Range("AC20").ClearContents()
$ws.Range("AE20").Value = $true
